$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("F2").Value = -5
$ws.Range("F5").Value = 0
$ws.Range("F11").Value = -4
$ws.Range("F15").Value = -2
$ws.Range("F16").Value = -8
$ws.Range("F17").Value = -4
$ws.Range("F19").Value = -2
$ws.Range("F21").Value = -1
$ws.Range("F23").Value = -8
$ws.Range("F24").Value = -4
$ws.Range("F26").Value = -4
$ws.Range("F28").Value = -6
$ws.Range("F29").Value = 2
$ws.Range("F35").Value = -4
$ws.Range("F36").Value = -4
$ws.Range("F38").Value = -3
$ws.Range("F41").Value = -3
$ws.Range("F43").Value = -4
$ws.Range("F44").Value = -2
$ws.Range("F45").Value = -7
$ws.Range("F46").Value = 11
$ws.Range("F47").Value = -4
$ws.Range("F48").Value = 1
$ws.Range("F52").Value = 9
$ws.Range("F57").Value = -4
$ws.Range("F61").Value = 0
$ws.Range("F64").Value = 5
$ws.Range("F67").Value = -8
